{"js": "// Insert a new \"Body Text\" paragraph right after the Introduction\n// paragraph (\"Git is a software tool used for Source Code Management...\")\n// and before the \"Terms\" heading, introducing that Git is a CLI tool.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the paragraph to anchor the insertion on, matched by its leading\n// text so the script is resilient to minor structural shifts.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"Git is a software tool used for Source Code Management\"\n    ) === 0\n  ) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the Introduction paragraph to anchor the new text on.\");\n}\n\n// Inserting \"After\" the target creates a new paragraph that inherits the\n// target's paragraph formatting (style \"Body Text\", left justification),\n// matching how Word splits a paragraph when you type a new one.\ntarget.insertParagraph(\n  \"Git is a CLI which stands for Command Line Interface.  Therefore, git needs to be run in a command prompt app or a terminal.  There are some GUI applications that can run git for you.  For instance, Visual Studio Code has a Git Extension for interacting with a git repository.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Introduction\" body-text paragraph that starts with\n# \"Git is a software tool used for Source Code Management\" and insert a\n# new Body Text paragraph right after it (before the \"Terms\" heading).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Git is a software tool used for Source Code Management*\") {\n        $target = $p\n        break\n    }\n}\n\n# Insert a new paragraph mark right after the target paragraph. The new\n# paragraph inherits the target's paragraph formatting (style \"Body Text\",\n# bidi/justification), matching how Word splits a paragraph.\n$target.Range.InsertParagraphAfter()\n\n# Find the newly created (still empty) paragraph immediately after the\n# target and give it the new sentence as its text.\n$idx = $target.Index\n$newPara = $d.Paragraphs($idx + 1)\n$newPara.Range.Text = \"Git is a CLI which stands for Command Line Interface.  Therefore, git needs to be run in a command prompt app or a terminal.  There are some GUI applications that can run git for you.  For instance, Visual Studio Code has a Git Extension for interacting with a git repository.\"\n"}
